$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 43.28265566666667
$ws.Range("N2").Value = 129.847967
$ws.Range("O2").Value = 0.667219228070094
$ws.Range("P2").Value = 0.667219228070094
$ws.Range("Q2").Value = 7.689769736362668
$ws.Range("R2").Value = 69.20792762726401
$ws.Range("S2").Value = 0.667219228070094
$ws.Range("T2").Value = 0.667219228070094

# Row 3
$ws.Range("O3").Value = 0.2872263480299067
$ws.Range("P3").Value = 0.2872263480299067
$ws.Range("S3").Value = 0.2872263480299067
$ws.Range("T3").Value = 0.2872263480299067

# Row 4
$ws.Range("M4").Value = 2.955125333333334
$ws.Range("N4").Value = 8.865376000000001
$ws.Range("O4").Value = 0.04555442389999943
$ws.Range("P4").Value = 0.04555442389999944
$ws.Range("Q4").Value = 0.5250193872213335
$ws.Range("R4").Value = 4.725174484992001
$ws.Range("S4").Value = 0.04555442389999943
$ws.Range("T4").Value = 0.04555442389999944
